$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Simplify the steel/industrial description in B2: remove the "RME" component
# from the "25% S/LFM+CDN/RME/H:1" line.
$text = "20% S+SL/LFM+CDN/H:1`n25% S/LFM+CDN/H:1`n25% CR+PC/LFM+CDN/H:1`n10% CR/LWAL+CDN/H:2`n8% CR/LFM+CDN/H:2`n5% W/LWAL+CDN/H:1`n7% MUR/LWAL+CDN/H:1"
$ws.Range("B2").Value = $text

# Wrap the (now multi-line) text and grow the row so it is fully visible.
$ws.Range("B2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 365

# Move the active selection to B10, matching the saved view state.
[void]$ws.Range("B10").Select()
